$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row: OBD2 extension cable ---
# Copy the formatting of the row above (C6) down to row 7 first so the
# new "Item" cell (C7) picks up the same style the other item cells use,
# then fill in the actual values.
$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)

$ws.Range("C7").Value = "OBD2 extension cable"
$ws.Range("D7").Value = 22.99

# --- New "Buyer" column (E) for every row ---
$ws.Range("E4").Value = "Buyer"
$ws.Range("E5").Value = "Mubashir"
$ws.Range("E6").Value = "Hussain"
$ws.Range("E7").Value = "Sarmad"
